# Changes of 31st March 2022
# Update the FedEx ShipmentTracking numbers (column P, rows 2-25) and the
# ActualRate (Q24) on Sheet1 to the new values. The tracking numbers are
# long digit strings that Excel would otherwise auto-convert to numbers,
# so the target ranges are forced to Text format before the values are
# written (matching how the source data is stored as shared strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New ShipmentTracking numbers for rows 2-25 (column P).
$trackingNumbers = @{
    2  = "320018191948"
    3  = "320018191959"
    4  = "320018191981"
    5  = "320018192006"
    6  = "320018192040"
    7  = "320018192061"
    8  = "320018192094"
    9  = "320018192131"
    10 = "320018192164"
    11 = "320018192186"
    12 = "320018192223"
    13 = "320018192245"
    14 = "320018192278"
    15 = "320018192290"
    16 = "320018192326"
    17 = "320018192348"
    18 = "320018192381"
    19 = "320018192407"
    20 = "320018192430"
    21 = "320018192451"
    22 = "320018192484"
    23 = "320018192495"
    24 = "320018192500"
    25 = "320018192510"
}

# Force column P (ShipmentTracking) on these rows to Text so the long
# numeric-looking tracking numbers are stored as strings, not numbers.
$ws.Range("P2:P25").NumberFormat = "@"

foreach ($row in $trackingNumbers.Keys) {
    $ws.Cells.Item($row, 16).Value = $trackingNumbers[$row]
}

# Update the ActualRate for row 24 as well. Keep it text (it's a
# formatted "$" amount stored as a shared string, not a number).
$ws.Range("Q24").NumberFormat = "@"
$ws.Range("Q24").Value = "$278.12"
